# Split the "Effective: ..." run so the year is wrapped in a
# gramStart/gramEnd proofing-error pair, matching the grammar checker's
# "comma before year" flag, e.g.:
#   Effective: October 31, [gramStart]2025[gramEnd] | v1.0
#
# Before:
#   <w:r><w:rPr><w:i/></w:rPr><w:t>Effective: October 31, 2025 | v1.0</w:t></w:r>
# After:
#   <w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Effective: October 31, </w:t></w:r>
#   <w:proofErr w:type="gramStart"/>
#   <w:r><w:rPr><w:i/></w:rPr><w:t>2025</w:t></w:r>
#   <w:proofErr w:type="gramEnd"/>
#   <w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> | v1.0</w:t></w:r>

$d = $word.ActiveDocument

# Locate the exact run text via Find (no replacement - just locate span).
$hit = $d.Content
$found = $hit.Find.Execute("Effective: October 31, 2025 | v1.0", $true, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Effective: ...' line to edit."
}

# Re-seat as a plain Range over the same span so InsertXML replaces
# (rather than appends after) the matched text.
$target = $d.Range($hit.Start, $hit.End)

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" ' + `
      'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Effective: October 31, </w:t></w:r>' + `
              '<w:proofErr w:type="gramStart"/>' + `
              '<w:r><w:rPr><w:i/></w:rPr><w:t>2025</w:t></w:r>' + `
              '<w:proofErr w:type="gramEnd"/>' + `
              '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> | v1.0</w:t></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$target.InsertXML($newXml)
